$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph (Paragraphs(1)).
# ------------------------------------------------------------------
$title = $d.Paragraphs(1)
$title.Range.InsertParagraphAfter()

# Grab the formatted text (empty run + bold run) from the paragraph
# that currently duplicates the title near the end of the document, so
# the new paragraph reuses the same run layout / bold formatting.
$countBefore = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($countBefore - 1)
$fmt = $dupTitlePara.Range.FormattedText

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaPara.Range.FormattedText = $fmt

# Swap the bold run's text for "Meta description"
$metaPara.Range.Find.Execute("Play Don Corlimone Free: Unique Fruit and Mafia-themed Slot", `
    $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

# Append the (non-bold) description text right after "Meta description",
# but before the paragraph mark.
$metaRange = $metaPara.Range
$insertPoint = $d.Range($metaRange.End - 1, $metaRange.End - 1)
$insertPoint.InsertAfter(": Explore the unique world of Don Corlimone, a mafia-themed slot with fruit-based customization. Play Don Corlimone free on desktop, mobile and tablet.")
$insertPoint.Bold = 0

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document, and rewrite the italic paragraph that follows it with
#    the new "Prompt: ..." image-generation description.
# ------------------------------------------------------------------
$countNow = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($countNow - 1)
$dupTitlePara.Range.Delete()

$countAfterDelete = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($countAfterDelete)
$italicRange = $italicPara.Range
$replaceRange = $d.Range($italicRange.Start, $italicRange.End - 1)
$replaceRange.Text = 'Prompt: Create a feature image for the game "Don Corlimone". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing next to a slot machine with mafia characters and fruit symbols on it. The background should be a cityscape reminiscent of New York, with tall buildings and the sky scrapers piercing the sky. The image should exude a fun and vibrant atmosphere that reflects the unique theme of the game.'
